# Sections no longer go across pages.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. edible_accessories: fill in real content (was a placeholder duplicate of
#    the monsters sheet) and widen column B so text wraps within the page.
# ---------------------------------------------------------------------------
$edible = $wb.Worksheets.Item("edible_accessories")
$edible.Range("A1:C2").ClearContents()
$edible.Range("A1").Value = "Accessory"
$edible.Range("B1").Value = "Effect"
$edible.Range("A2").Value = "Myrrh Leaf"
$edible.Range("B2").Value = "Clear a Status Effect from yourself"
$edible.Columns.Item(2).ColumnWidth = 29.53

# ---------------------------------------------------------------------------
# 2. monsters: replace the placeholder content with the real monster stat
#    block (Goblin).
# ---------------------------------------------------------------------------
$monsters = $wb.Worksheets.Item("monsters")
$monsters.Range("A1:C2").ClearContents()
$monsters.Range("A1").Value = "Monster"
$monsters.Range("B1").Value = "Hearts"
$monsters.Range("C1").Value = "Combat"
$monsters.Range("D1").Value = "Unity"
$monsters.Range("E1").Value = "Magic"
$monsters.Range("F1").Value = "Resistances"
$monsters.Range("G1").Value = "Vulnerabilities"
$monsters.Range("H1").Value = "Abilities"
$monsters.Range("A2").Value = "Goblin"
$monsters.Range("B2").Value = 3
$monsters.Range("C2").Value = 2
$monsters.Range("D2").Value = 1
$monsters.Range("E2").Value = 1
$monsters.Columns.Item(7).ColumnWidth = 12.69
$monsters.Columns.Item(8).ColumnWidth = 23.41

# ---------------------------------------------------------------------------
# 3. New sheet "monsters_rva" (resistances/vulnerabilities/abilities detail),
#    placed right after "monsters". Add + move first, fill in data after, so
#    sheet-collection indices used by Move() are never stale.
# ---------------------------------------------------------------------------
$rva = $wb.Worksheets.Add()
$rva.Name = "monsters_rva"
$rva.Move($null, $wb.Worksheets.Item("monsters"))

# ---------------------------------------------------------------------------
# 4. New sheet "Sheet10" with an imported table (Import1) plus a stray cell,
#    placed right after "monsters_rva" and left as the active sheet/tab.
# ---------------------------------------------------------------------------
$s10 = $wb.Worksheets.Add()
$s10.Name = "Sheet10"
$s10.Move($null, $wb.Worksheets.Item("monsters_rva"))

$rva = $wb.Worksheets.Item("monsters_rva")
$rva.Range("A1").Value = "Monster"
$rva.Range("B1").Value = "Details"
$rva.Range("C1").Value = "Type"
$rva.Range("A2").Value = "Goblin"
$rva.Range("B2").Value = "Fire"
$rva.Range("C2").Value = "Vulnerable"
$rva.Range("A3").Value = "Goblin"
$rva.Range("B3").Value = "Lightning"
$rva.Range("C3").Value = "Vulnerable"
$rva.Range("A4").Value = "Goblin"
$rva.Range("B4").Value = "Appears in groups of 3 to 4"
$rva.Range("C4").Value = "Ability"

$rva.Columns.Item(2).ColumnWidth = 23.41
$rva.Columns.Item(7).ColumnWidth = 12.69
$rva.Columns.Item(8).ColumnWidth = 23.41

$s10 = $wb.Worksheets.Item("Sheet10")
$s10.Range("A1").Value = "ID"
$s10.Range("B1").Value = "Type"
$s10.Range("C1").Value = "Description"
$s10.Range("D1").Value = "Monster"
$s10.Range("A2").Value = 1
$s10.Range("B2").Value = "Goblin"
$s10.Range("C2").Value = "Fire"
$s10.Range("D2").Value = "Vulnerable"
$s10.Range("B3").Value = "asdf"

$lo = $s10.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $s10.Range("A1:D2"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Import1"

$s10.Range("A1").Select()
$s10.Activate()

Write-Host "done"
